# "Pais" worksheet refresh: newer COVID-19 snapshot timestamp + per-country stats,
# including a handful of countries whose ranking (row position) changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp caption in A1
$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 1 de Agosto de 2020 a las 01:25'

# Country rows: label (if the country occupying this row changed) + the 7 stat columns
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 4703743
$ws.Cells.Item(4, 3).Value = 68758
$ws.Cells.Item(4, 4).Value = 2325306
$ws.Cells.Item(4, 5).Value = 2221702
$ws.Cells.Item(4, 7).Value = 1450
$ws.Cells.Item(4, 8).Value = 156735

# Row 5: Brasil
$ws.Cells.Item(5, 2).Value = 2666298
$ws.Cells.Item(5, 3).Value = 52509
$ws.Cells.Item(5, 4).Value = 1884051
$ws.Cells.Item(5, 5).Value = 689679
$ws.Cells.Item(5, 7).Value = 1191
$ws.Cells.Item(5, 8).Value = 92568

# Row 6: India
$ws.Cells.Item(6, 2).Value = 1697054
$ws.Cells.Item(6, 3).Value = 57704
$ws.Cells.Item(6, 5).Value = 564856

# Row 22: Argentina
$ws.Cells.Item(22, 1).Value = 'Argentina'
$ws.Cells.Item(22, 2).Value = 191302
$ws.Cells.Item(22, 3).Value = 5929
$ws.Cells.Item(22, 4).Value = 83780
$ws.Cells.Item(22, 5).Value = 103979
$ws.Cells.Item(22, 7).Value = 102
$ws.Cells.Item(22, 8).Value = 3543

# Row 23: Francia
$ws.Cells.Item(23, 1).Value = 'Francia'
$ws.Cells.Item(23, 2).Value = 187919
$ws.Cells.Item(23, 3).Value = 1346
$ws.Cells.Item(23, 4).Value = 81500
$ws.Cells.Item(23, 5).Value = 76154
$ws.Cells.Item(23, 7).Value = 11
$ws.Cells.Item(23, 8).Value = 30265

# Row 42: Panama
$ws.Cells.Item(42, 2).Value = 65256
$ws.Cells.Item(42, 3).Value = 1065
$ws.Cells.Item(42, 4).Value = 39166
$ws.Cells.Item(42, 5).Value = 24669
$ws.Cells.Item(42, 7).Value = 24
$ws.Cells.Item(42, 8).Value = 1421

# Row 50: Nigeria
$ws.Cells.Item(50, 2).Value = 43151
$ws.Cells.Item(50, 3).Value = 462
$ws.Cells.Item(50, 4).Value = 19565
$ws.Cells.Item(50, 5).Value = 22707
$ws.Cells.Item(50, 7).Value = 1
$ws.Cells.Item(50, 8).Value = 879

# Row 58: Japon
$ws.Cells.Item(58, 2).Value = 34372
$ws.Cells.Item(58, 3).Value = 1323
$ws.Cells.Item(58, 4).Value = 24929
$ws.Cells.Item(58, 5).Value = 8437
$ws.Cells.Item(58, 7).Value = 2
$ws.Cells.Item(58, 8).Value = 1006

# Row 69: Venezuela
$ws.Cells.Item(69, 2).Value = 18574
$ws.Cells.Item(69, 3).Value = 715
$ws.Cells.Item(69, 5).Value = 7989
$ws.Cells.Item(69, 7).Value = 6
$ws.Cells.Item(69, 8).Value = 164

# Row 75: Chequia
$ws.Cells.Item(75, 2).Value = 16574
$ws.Cells.Item(75, 3).Value = 203
$ws.Cells.Item(75, 4).Value = 11569
$ws.Cells.Item(75, 5).Value = 4623

# Row 81: Sudan
$ws.Cells.Item(81, 2).Value = 11644
$ws.Cells.Item(81, 3).Value = 148
$ws.Cells.Item(81, 4).Value = 6119
$ws.Cells.Item(81, 5).Value = 4779
$ws.Cells.Item(81, 7).Value = 21
$ws.Cells.Item(81, 8).Value = 746

# Row 86: Noruega
$ws.Cells.Item(86, 2).Value = 9240
$ws.Cells.Item(86, 3).Value = 32
$ws.Cells.Item(86, 5).Value = 233

# Row 89: Guayana Francesa
$ws.Cells.Item(89, 2).Value = 7799
$ws.Cells.Item(89, 3).Value = 71
$ws.Cells.Item(89, 4).Value = 6386
$ws.Cells.Item(89, 5).Value = 1370

# Row 102: Guinea Ecuatorial
$ws.Cells.Item(102, 1).Value = 'Guinea Ecuatorial'
$ws.Cells.Item(102, 2).Value = 4821
$ws.Cells.Item(102, 3).Value = 1750
$ws.Cells.Item(102, 4).Value = 2182
$ws.Cells.Item(102, 5).Value = 2556
$ws.Cells.Item(102, 7).Value = 32
$ws.Cells.Item(102, 8).Value = 83

# Row 103: Republica de Africa Central
$ws.Cells.Item(103, 1).Value = 'Republica de Africa Central'
$ws.Cells.Item(103, 2).Value = 4608
$ws.Cells.Item(103, 3).Value = 3
$ws.Cells.Item(103, 4).Value = 1606
$ws.Cells.Item(103, 5).Value = 2943
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 59

# Row 104: Libano
$ws.Cells.Item(104, 1).Value = 'Libano'
$ws.Cells.Item(104, 2).Value = 4555
$ws.Cells.Item(104, 3).Value = 221
$ws.Cells.Item(104, 4).Value = 1761
$ws.Cells.Item(104, 5).Value = 2733
$ws.Cells.Item(104, 7).Value = 4
$ws.Cells.Item(104, 8).Value = 61

# Row 105: Hungria
$ws.Cells.Item(105, 1).Value = 'Hungria'
$ws.Cells.Item(105, 2).Value = 4505
$ws.Cells.Item(105, 3).Value = 21
$ws.Cells.Item(105, 4).Value = 3353
$ws.Cells.Item(105, 5).Value = 556
$ws.Cells.Item(105, 7).Value = 0
$ws.Cells.Item(105, 8).Value = 596

# Row 106: Grecia
$ws.Cells.Item(106, 1).Value = 'Grecia'
$ws.Cells.Item(106, 2).Value = 4477
$ws.Cells.Item(106, 3).Value = 76
$ws.Cells.Item(106, 4).Value = 1374
$ws.Cells.Item(106, 5).Value = 2897
$ws.Cells.Item(106, 7).Value = 3
$ws.Cells.Item(106, 8).Value = 206

# Row 107: Malaui
$ws.Cells.Item(107, 1).Value = 'Malaui'
$ws.Cells.Item(107, 2).Value = 4078
$ws.Cells.Item(107, 3).Value = 220
$ws.Cells.Item(107, 4).Value = 1875
$ws.Cells.Item(107, 5).Value = 2089
$ws.Cells.Item(107, 7).Value = 7
$ws.Cells.Item(107, 8).Value = 114

# Row 108: Maldivas
$ws.Cells.Item(108, 1).Value = 'Maldivas'
$ws.Cells.Item(108, 2).Value = 3793
$ws.Cells.Item(108, 3).Value = 74
$ws.Cells.Item(108, 4).Value = 2607
$ws.Cells.Item(108, 5).Value = 1170
$ws.Cells.Item(108, 8).Value = 16

# Row 109: Nicaragua
$ws.Cells.Item(109, 1).Value = 'Nicaragua'
$ws.Cells.Item(109, 2).Value = 3672
$ws.Cells.Item(109, 3).Value = 0
$ws.Cells.Item(109, 4).Value = 2492
$ws.Cells.Item(109, 5).Value = 1064
$ws.Cells.Item(109, 7).Value = 0
$ws.Cells.Item(109, 8).Value = 116

# Row 110: Libia
$ws.Cells.Item(110, 1).Value = 'Libia'
$ws.Cells.Item(110, 2).Value = 3621
$ws.Cells.Item(110, 3).Value = 183
$ws.Cells.Item(110, 4).Value = 618
$ws.Cells.Item(110, 5).Value = 2929
$ws.Cells.Item(110, 7).Value = 1
$ws.Cells.Item(110, 8).Value = 74

# Row 111: Tailandia
$ws.Cells.Item(111, 1).Value = 'Tailandia'
$ws.Cells.Item(111, 2).Value = 3310
$ws.Cells.Item(111, 3).Value = 6
$ws.Cells.Item(111, 4).Value = 3125
$ws.Cells.Item(111, 5).Value = 127
$ws.Cells.Item(111, 7).Value = 0
$ws.Cells.Item(111, 8).Value = 58

# Row 112: Hong Kong
$ws.Cells.Item(112, 1).Value = 'Hong Kong'
$ws.Cells.Item(112, 2).Value = 3273
$ws.Cells.Item(112, 3).Value = 121
$ws.Cells.Item(112, 4).Value = 1751
$ws.Cells.Item(112, 5).Value = 1495
$ws.Cells.Item(112, 7).Value = 2
$ws.Cells.Item(112, 8).Value = 27

# Row 113: Somalia
$ws.Cells.Item(113, 1).Value = 'Somalia'
$ws.Cells.Item(113, 2).Value = 3212
$ws.Cells.Item(113, 4).Value = 1562
$ws.Cells.Item(113, 5).Value = 1557
$ws.Cells.Item(113, 8).Value = 93

# Row 114: Congo
$ws.Cells.Item(114, 1).Value = 'Congo'
$ws.Cells.Item(114, 2).Value = 3200
$ws.Cells.Item(114, 4).Value = 829
$ws.Cells.Item(114, 5).Value = 2317
$ws.Cells.Item(114, 8).Value = 54

# Row 115: Zimbabue
$ws.Cells.Item(115, 1).Value = 'Zimbabue'
$ws.Cells.Item(115, 2).Value = 3169
$ws.Cells.Item(115, 3).Value = 77
$ws.Cells.Item(115, 4).Value = 1004
$ws.Cells.Item(115, 5).Value = 2098
$ws.Cells.Item(115, 7).Value = 14
$ws.Cells.Item(115, 8).Value = 67

# Row 116: Montenegro
$ws.Cells.Item(116, 1).Value = 'Montenegro'
$ws.Cells.Item(116, 2).Value = 3073
$ws.Cells.Item(116, 3).Value = 57
$ws.Cells.Item(116, 4).Value = 1005
$ws.Cells.Item(116, 5).Value = 2020
$ws.Cells.Item(116, 7).Value = 1
$ws.Cells.Item(116, 8).Value = 48

# Row 122: Cabo Verde
$ws.Cells.Item(122, 4).Value = 1824
$ws.Cells.Item(122, 5).Value = 604

# Row 145: Angola
$ws.Cells.Item(145, 1).Value = 'Angola'
$ws.Cells.Item(145, 2).Value = 1148
$ws.Cells.Item(145, 3).Value = 39
$ws.Cells.Item(145, 4).Value = 437
$ws.Cells.Item(145, 5).Value = 659
$ws.Cells.Item(145, 7).Value = 1
$ws.Cells.Item(145, 8).Value = 52

# Row 146: Niger
$ws.Cells.Item(146, 1).Value = 'Niger'
$ws.Cells.Item(146, 2).Value = 1134
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 1028
$ws.Cells.Item(146, 5).Value = 37
$ws.Cells.Item(146, 8).Value = 69

# Row 147: Republica de Chipre
$ws.Cells.Item(147, 1).Value = 'Republica de Chipre'
$ws.Cells.Item(147, 2).Value = 1114
$ws.Cells.Item(147, 3).Value = 24
$ws.Cells.Item(147, 4).Value = 852
$ws.Cells.Item(147, 5).Value = 243
$ws.Cells.Item(147, 8).Value = 19

# Row 149: Togo
$ws.Cells.Item(149, 1).Value = 'Togo'
$ws.Cells.Item(149, 2).Value = 941
$ws.Cells.Item(149, 3).Value = 33
$ws.Cells.Item(149, 4).Value = 641
$ws.Cells.Item(149, 5).Value = 281
$ws.Cells.Item(149, 7).Value = 1
$ws.Cells.Item(149, 8).Value = 19

# Row 150: Republica del Chad
$ws.Cells.Item(150, 1).Value = 'Republica del Chad'
$ws.Cells.Item(150, 2).Value = 936
$ws.Cells.Item(150, 3).Value = 1
$ws.Cells.Item(150, 4).Value = 813
$ws.Cells.Item(150, 5).Value = 48
$ws.Cells.Item(150, 8).Value = 75

# Row 152: Santo Tome y Principe
$ws.Cells.Item(152, 2).Value = 871
$ws.Cells.Item(152, 3).Value = 1
$ws.Cells.Item(152, 4).Value = 778
$ws.Cells.Item(152, 5).Value = 78

# Row 161: Bahamas
$ws.Cells.Item(161, 1).Value = 'Bahamas'
$ws.Cells.Item(161, 2).Value = 574
$ws.Cells.Item(161, 3).Value = 66
$ws.Cells.Item(161, 4).Value = 91
$ws.Cells.Item(161, 5).Value = 469
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 14

# Row 162: Vietnam
$ws.Cells.Item(162, 1).Value = 'Vietnam'
$ws.Cells.Item(162, 2).Value = 546
$ws.Cells.Item(162, 3).Value = 37
$ws.Cells.Item(162, 4).Value = 373
$ws.Cells.Item(162, 5).Value = 171
$ws.Cells.Item(162, 7).Value = 2
$ws.Cells.Item(162, 8).Value = 2

# Row 163: Tanzania
$ws.Cells.Item(163, 1).Value = 'Tanzania'
$ws.Cells.Item(163, 2).Value = 509
$ws.Cells.Item(163, 4).Value = 183
$ws.Cells.Item(163, 5).Value = 305
$ws.Cells.Item(163, 8).Value = 21

# Row 164: Gambia
$ws.Cells.Item(164, 1).Value = 'Gambia'
$ws.Cells.Item(164, 2).Value = 498
$ws.Cells.Item(164, 3).Value = 95
$ws.Cells.Item(164, 4).Value = 68
$ws.Cells.Item(164, 5).Value = 421
$ws.Cells.Item(164, 7).Value = 1
$ws.Cells.Item(164, 8).Value = 9

# Row 165: Taiwan
$ws.Cells.Item(165, 1).Value = 'Taiwan'
$ws.Cells.Item(165, 2).Value = 467
$ws.Cells.Item(165, 4).Value = 441
$ws.Cells.Item(165, 5).Value = 19
$ws.Cells.Item(165, 8).Value = 7

